# refactoring for N vs N mode
# - sheet1 (chesstactic_soldier): "side" column switches from text labels
#   ("my"/"opp") to a numeric side id (0 = my team, 1 = opponent team).
# - sheet2 (chesstactic_tactic): same numeric "side" id, plus a new "name"
#   column (team display name) inserted between "side" and "attack".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: chesstactic_soldier -----------------------------------------
# Column A (side) rows 2-7 were "my" -> now numeric 0
for ($r = 2; $r -le 7; $r++) {
  $ws1.Cells.Item($r, 1).Value = 0
}
# Column A (side) rows 8-13 were "opp" -> now numeric 1
for ($r = 8; $r -le 13; $r++) {
  $ws1.Cells.Item($r, 1).Value = 1
}

# --- Sheet2: chesstactic_tactic -------------------------------------------
# Insert a new "name" column between "side" (A) and "attack" (B)
$null = $ws2.Columns("B:B").Insert()

$ws2.Cells.Item(1, 1).Value = "side"
$ws2.Cells.Item(1, 2).Value = "name"
$ws2.Cells.Item(1, 3).Value = "attack"
$ws2.Cells.Item(1, 4).Value = "defence"

$ws2.Cells.Item(2, 1).Value = 0
$ws2.Cells.Item(2, 2).Value = "내팀"
$ws2.Cells.Item(2, 3).Value = 0
$ws2.Cells.Item(2, 4).Value = 0

$ws2.Cells.Item(3, 1).Value = 1
$ws2.Cells.Item(3, 2).Value = "적1"
$ws2.Cells.Item(3, 3).Value = 0
$ws2.Cells.Item(3, 4).Value = 0

# --- View state -------------------------------------------------------------
# Sheet2 selection moves to A11; sheet1 stays the active tab with
# its frozen header row and selection on G13.
$null = $ws2.Range("A11").Select()
$null = $ws1.Activate()
$null = $ws1.Range("G13").Select()

Write-Output "chesstactic config updated for N vs N mode"
